$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P16").Value = 7.26666666666667
$ws.Range("P17").Value = 6.66666666666667
$ws.Range("P18").Value = 6.26666666666667
$ws.Range("P20").Value = 5.43333333333333
$ws.Range("P24").Value = 4.23333333333333
$ws.Range("P26").Value = 3.96666666666667
$ws.Range("P28").Value = 4.36666666666667
$ws.Range("P30").Value = 4.53333333333333
$ws.Range("P31").Value = 4.63333333333333
$ws.Range("P33").Value = 5.26666666666667
$ws.Range("P34").Value = 4.83333333333333
$ws.Range("P35").Value = 4.76666666666667
$ws.Range("P37").Value = 4.73333333333333
$ws.Range("P38").Value = 5.13333333333333
$ws.Range("P39").Value = 4.66666666666667
$ws.Range("P42").Value = 4.43333333333333
$ws.Range("P45").Value = 4.53333333333333
$ws.Range("P47").Value = 4.73333333333333
$ws.Range("P49").Value = 4.63333333333333
$ws.Range("P53").Value = 5.23333333333333
$ws.Range("P54").Value = 5.36666666666667
$ws.Range("P55").Value = 5.93333333333333
$ws.Range("P56").Value = 7.36666666666667
$ws.Range("P57").Value = 8.53333333333333
$ws.Range("P58").Value = 10.8666666666667
$ws.Range("P61").Value = 13.8333333333333
$ws.Range("P63").Value = 14.3666666666667
$ws.Range("P65").Value = 15.6333333333333
$ws.Range("P66").Value = 15.1333333333333
$ws.Range("P67").Value = 15.0333333333333
$ws.Range("P69").Value = 15.8333333333333
$ws.Range("I70").Value = 22743
$ws.Range("P70").Value = 15.9333333333333
$ws.Range("I71").Value = 22944
$ws.Range("P71").Value = 15.5333333333333
$ws.Range("I72").Value = 22771
$ws.Range("P72").Value = 15.4666666666667
$ws.Range("I73").Value = 22182
$ws.Range("I74").Value = 22435
$ws.Range("P74").Value = 14.5666666666667
$ws.Range("I75").Value = 23226
$ws.Range("P75").Value = 14.3666666666667
$ws.Range("I76").Value = 22646
$ws.Range("P76").Value = 13.3333333333333
$ws.Range("I77").Value = 22431
$ws.Range("P77").Value = 12.8333333333333
$ws.Range("I78").Value = 22588
$ws.Range("P78").Value = 12.8333333333333
$ws.Range("I79").Value = 23272
$ws.Range("I80").Value = 23061
$ws.Range("P80").Value = 11.6666666666667
$ws.Range("I81").Value = 22948
$ws.Range("P81").Value = 10.8666666666667
$ws.Range("I82").Value = 23660
$ws.Range("P82").Value = 10.6333333333333
$ws.Range("I83").Value = 24848
$ws.Range("I84").Value = 24816
$ws.Range("P84").Value = 9.53333333333333
$ws.Range("I85").Value = 23665
$ws.Range("P85").Value = 9.43333333333333
$ws.Range("I86").Value = 24772
$ws.Range("P86").Value = 8.93333333333333
$ws.Range("I87").Value = 25687
$ws.Range("P87").Value = 8.86666666666667
$ws.Range("I88").Value = 26358
$ws.Range("P88").Value = 8.23333333333333
$ws.Range("I89").Value = 24939
$ws.Range("P89").Value = 7.56666666666667
$ws.Range("I90").Value = 26395
$ws.Range("P90").Value = 7.26666666666667
$ws.Range("I91").Value = 27309
$ws.Range("P91").Value = 6.66666666666667
$ws.Range("I92").Value = 27844
$ws.Range("I93").Value = 26734
$ws.Range("P93").Value = 6.43333333333333
$ws.Range("I94").Value = 27081
$ws.Range("P94").Value = 5.93333333333333
$ws.Range("I95").Value = 28748
$ws.Range("P95").Value = 5.86666666666667
$ws.Range("I96").Value = 29019
$ws.Range("I97").Value = 28081
$ws.Range("P97").Value = 5.73333333333333
$ws.Range("I98").Value = 29274
$ws.Range("P98").Value = 4.96666666666667
$ws.Range("I99").Value = 30806
$ws.Range("P99").Value = 5.16666666666667
$ws.Range("I100").Value = 30628
$ws.Range("I101").Value = 30185
$ws.Range("I102").Value = 32262
$ws.Range("P102").Value = 4.86666666666667
$ws.Range("I103").Value = 31972
$ws.Range("M103").Value = 33.3
$ws.Range("P103").Value = 5.13333333333333
$ws.Range("I104").Value = 31069
$ws.Range("I105").Value = 30806
$ws.Range("P105").Value = 6.26666666666667
$ws.Range("I106").Value = 31237
$ws.Range("P106").Value = 7.33333333333333
$ws.Range("I107").Value = 33802
$ws.Range("P107").Value = 6.96666666666667
$ws.Range("I108").Value = 34570
$ws.Range("P108").Value = 5.46666666666667
$ws.Range("I109").Value = 33015
$ws.Range("P109").Value = 5.23333333333333
$ws.Range("I110").Value = 34043
$ws.Range("M110").Value = 15.8
$ws.Range("P110").Value = 4.86666666666667
$ws.Range("I111").Value = 36627
$ws.Range("M111").Value = 15.8
$ws.Range("I112").Value = 37065
$ws.Range("P112").Value = 4.26666666666667
$ws.Range("I113").Value = 35289
$ws.Range("P113").Value = 4.46666666666667
$ws.Range("I114").Value = 38607
$ws.Range("I115").Value = 40465
$ws.Range("M115").Value = 11.5
$ws.Range("P115").Value = 4.16666666666667
$ws.Range("I116").Value = 40861
$ws.Range("M116").Value = 12.2
$ws.Range("P116").Value = 4.43333333333333
$ws.Range("F117").Value = 124.8
$ws.Range("I117").Value = 39549
$ws.Range("M117").Value = 11
$ws.Range("P117").Value = 4.43333333333333
$ws.Range("I118").Value = 43059
$ws.Range("M118").Value = 14.6
$ws.Range("P118").Value = 4.26666666666667
$ws.Range("T118").Value = 53998355
$ws.Range("I119").Value = 43618
$ws.Range("T119").Value = 54218647
$ws.Range("I120").Value = 44381
$ws.Range("P120").Value = 4.23333333333333
$ws.Range("T120").Value = 57909652
$ws.Range("F121").Value = 141.2
$ws.Range("I121").Value = 43133
$ws.Range("M121").Value = 13.4
$ws.Range("P121").Value = 4.26666666666667
$ws.Range("T121").Value = 57607313
$ws.Range("F122").Value = 140.7
$ws.Range("I122").Value = 44847
$ws.Range("M122").Value = 13.4
$ws.Range("P122").Value = 4.43333333333333
$ws.Range("S122").Value = 1201.926
$ws.Range("T122").Value = 87996779
$ws.Range("F123").Value = 150.3
$ws.Range("I123").Value = 45928
$ws.Range("O123").Value = 2554266.66666667
$ws.Range("L124").Value = 122.533333333333
$ws.Range("P124").Value = 4.73333333333333
